# This change performs a cyclic shift (rotation) of the per-observation
# data in rows 3-6: row 3's data moves to row 4, row 4's to row 5,
# row 5's to row 6, and row 6's data wraps back around to row 3.
# Only the observation-specific columns move (A, B, E, F, G, H, Q, R, AC);
# all other columns (shared metadata like locality, date, observer, etc.)
# are identical across these rows and therefore stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 4, 5, 6)
$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "AC")

# Capture the current ("before") values for every relevant cell so we can
# write them back out in rotated order without retyping/transcribing them.
$snapshot = @{}
foreach ($r in $rows) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $addr = "$col$r"
        $rowVals[$col] = $ws.Range($addr).Value2
    }
    $snapshot[$r] = $rowVals
}

# Write each row's values into the next row down, with row 6 wrapping
# around to row 3.
for ($i = 0; $i -lt $rows.Length; $i++) {
    $destRow = $rows[$i]
    $srcRow = $rows[($i - 1 + $rows.Length) % $rows.Length]
    $srcVals = $snapshot[$srcRow]

    foreach ($col in $cols) {
        $addr = "$col$destRow"
        $val = $srcVals[$col]
        if ($col -eq "AC" -and ($null -eq $val -or $val -eq "")) {
            $ws.Range($addr).ClearContents()
        } else {
            $ws.Range($addr).Value = $val
        }
    }
}
